# Update the "想去人数" (want-to-go count) figures that changed between
# the two gh-pages data refreshes.
#
# Sheet "展览"    (Exhibition) rows 2,4,5,8 -> column F
# Sheet "全部类型" (All types)  rows 2,4,5,9 -> column F
# (the "全部类型" sheet has one extra row coming from the "演出" sheet,
#  which is why the row numbers are shifted by one there)

$wb = $excel.ActiveWorkbook

function Update-Count($ws, $row, $oldVal, $newVal) {
    $cell = $ws.Cells.Item($row, 6)
    $current = $cell.Value()
    if ($current -eq $oldVal) {
        $cell.Value = $newVal
    }
}

$ws1 = $wb.Worksheets.Item("展览")
Update-Count $ws1 2 293  294
Update-Count $ws1 4 2431 2442
Update-Count $ws1 5 1785 1789
Update-Count $ws1 8 840  845

$ws4 = $wb.Worksheets.Item("全部类型")
Update-Count $ws4 2 293  294
Update-Count $ws4 4 2431 2442
Update-Count $ws4 5 1785 1789
Update-Count $ws4 9 840  845
